$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999998603162332,
    0.9951161208588728,
    0.97511854951286581,
    0.96623165583902104,
    0.95778689508192838,
    0.94001098753227108,
    0.93841893926472908,
    0.93428345976429417,
    0.93234449904993388,
    0.93133580039802566,
    0.93109830692008499,
    0.93090315067861285,
    0.93019897710488963,
    0.93048891967413527,
    0.93225326839256151,
    0.93500954435269767,
    0.94344574959887217,
    0.94916297287933049,
    0.99309850810569578,
    0.97736203349048323,
    0.97596358544273865,
    0.97469907850577886,
    0.96808324515000854,
    0.95506182172925191,
    0.94860475749632234,
    0.93635563677796729,
    0.93150909111297542,
    0.91003625398423016,
    0.89476520434479867,
    0.88819498955918263,
    0.88054160950272753,
    0.8788623792632575,
    0.87834239579793194
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
